# Applies the "fixed workflow" re-run: for both the NBR and BAR sheets,
# the Cutoff column (B) and Reaction_number column (C) data now come
# from what used to be rows 6-20 (i.e. shifted up by 4 rows), and the
# trailing 4 rows (17-20) are removed, shrinking the sheet from
# A1:C20 down to A1:C16. Column A (Threshold index) is left as-is
# (0..14) since it is just a 0-based row counter.

$wb = $excel.ActiveWorkbook

# New B (Cutoff) / C (Reaction_number) values for rows 2..16 on each sheet,
# taken directly from what used to be rows 6..20 before the fix.
$nbrValues = @(
    @(5, 105),
    @(6, 103),
    @(7, 103),
    @(8, 102),
    @(9, 102),
    @(10, 102),
    @(11, 101),
    @(12, 99),
    @(13, 97),
    @(14, 96),
    @(15, 96),
    @(16, 96),
    @(17, 95),
    @(18, 95),
    @(19, 95)
)

$barValues = @(
    @(5, 592),
    @(6, 595),
    @(7, 595),
    @(8, 594),
    @(9, 594),
    @(10, 588),
    @(11, 589),
    @(12, 589),
    @(13, 591),
    @(14, 589),
    @(15, 588),
    @(16, 588),
    @(17, 587),
    @(18, 589),
    @(19, 590)
)

foreach ($sheetInfo in @(
        @{ Name = "NBR"; Values = $nbrValues },
        @{ Name = "BAR"; Values = $barValues }
    )) {

    $ws = $wb.Worksheets.Item($sheetInfo.Name)
    $values = $sheetInfo.Values

    # Write the updated B/C values into rows 2..16.
    for ($i = 0; $i -lt $values.Count; $i++) {
        $row = 2 + $i
        $pair = $values[$i]
        $ws.Cells.Item($row, 2).Value = $pair[0]
        $ws.Cells.Item($row, 3).Value = $pair[1]
    }

    # Remove the now-stale trailing rows 17..20.
    $ws.Rows("17:20").Delete()
}
